$wb = $excel.ActiveWorkbook

# --- Step 1: insert the new "test_suite" sheet before the first sheet ---
# NOTE: worksheet variables captured before Add() become positionally stale
# (they track "whatever sheet sits at that slot") once a new sheet is
# spliced in, so every other sheet is re-fetched by name on demand below
# instead of being cached in a variable.
$wsSuite = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$wsSuite.Name = "test_suite"

# --- Step 2: write cells in the precise order needed so new shared strings ---
# --- land at the same indices as the target workbook.                   ---

# OpenAccountTest: fix existing "rupee" -> "Rupee" (modifies shared string in place)
$wb.Worksheets.Item("OpenAccountTest").Range("B2").Value = "Rupee"

# New shared strings 12,13 (OpenAccountTest header/value for alertMessage column)
$wb.Worksheets.Item("OpenAccountTest").Range("C2").Value = "Account created successfully"
$wb.Worksheets.Item("OpenAccountTest").Range("C1").Value = "alertMessage"

# New shared strings 14,15,16 (AddCustomerTest new rows, first name column)
$wb.Worksheets.Item("AddCustomerTest").Range("A3").Value = "rakesh"
$wb.Worksheets.Item("AddCustomerTest").Range("A4").Value = "sudesh"
$wb.Worksheets.Item("AddCustomerTest").Range("A5").Value = "naveen"

# New shared strings 17,18,19 (OpenAccountTest new rows, customer column)
$wb.Worksheets.Item("OpenAccountTest").Range("A3").Value = "rakesh kumar"
$wb.Worksheets.Item("OpenAccountTest").Range("A4").Value = "sudesh kumar"
$wb.Worksheets.Item("OpenAccountTest").Range("A5").Value = "naveen kumar"

# New shared strings 20,21,22,23,24,25 (test_suite sheet content)
# Column A is populated first in full (TCID/BankManagerLoginTest/AddCustomerTest/
# OpenAccountTest) so those four strings claim indices 20,22,23,24 - with
# "runMode" (21) wedged in via B1 - before "Y" (25) is minted from B2.
$wsSuite.Range("A1").Value = "TCID"
$wsSuite.Range("B1").Value = "runMode"
$wsSuite.Range("A2").Value = "BankManagerLoginTest"
$wsSuite.Range("A3").Value = "AddCustomerTest"
$wsSuite.Range("A4").Value = "OpenAccountTest"
$wsSuite.Range("B2").Value = "Y"
$wsSuite.Range("B3").Value = "Y"
$wsSuite.Range("B4").Value = "Y"

# New shared string 26 ("N") - first used on AddCustomerTest's runMode column
$wb.Worksheets.Item("AddCustomerTest").Range("E3").Value = "N"

# --- Step 3: fill in the remaining cells, reusing shared strings created above ---

# AddCustomerTest: finish rows 3-5 (existing columns B-D) + new runMode column E
$wb.Worksheets.Item("AddCustomerTest").Range("B3").Value = "kumar"
$wb.Worksheets.Item("AddCustomerTest").Range("C3").Value = "sfwr34d"
$wb.Worksheets.Item("AddCustomerTest").Range("D3").Value = "Customer added successfully"
$wb.Worksheets.Item("AddCustomerTest").Range("B4").Value = "kumar"
$wb.Worksheets.Item("AddCustomerTest").Range("C4").Value = "sfwr34d"
$wb.Worksheets.Item("AddCustomerTest").Range("D4").Value = "Customer added successfully"
$wb.Worksheets.Item("AddCustomerTest").Range("B5").Value = "kumar"
$wb.Worksheets.Item("AddCustomerTest").Range("C5").Value = "sfwr34d"
$wb.Worksheets.Item("AddCustomerTest").Range("D5").Value = "Customer added successfully"

$wb.Worksheets.Item("AddCustomerTest").Range("E1").Value = "runMode"
$wb.Worksheets.Item("AddCustomerTest").Range("E2").Value = "Y"
$wb.Worksheets.Item("AddCustomerTest").Range("E4").Value = "Y"
$wb.Worksheets.Item("AddCustomerTest").Range("E5").Value = "Y"

# OpenAccountTest: finish rows 3-5 (existing columns A,B) + new alertMessage/runMode columns
$wb.Worksheets.Item("OpenAccountTest").Range("B3").Value = "Rupee"
$wb.Worksheets.Item("OpenAccountTest").Range("C3").Value = "Account created successfully"
$wb.Worksheets.Item("OpenAccountTest").Range("B4").Value = "Rupee"
$wb.Worksheets.Item("OpenAccountTest").Range("C4").Value = "Account created successfully"
$wb.Worksheets.Item("OpenAccountTest").Range("B5").Value = "Rupee"
$wb.Worksheets.Item("OpenAccountTest").Range("C5").Value = "Account created successfully"

$wb.Worksheets.Item("OpenAccountTest").Range("D1").Value = "runMode"
$wb.Worksheets.Item("OpenAccountTest").Range("D2").Value = "Y"
$wb.Worksheets.Item("OpenAccountTest").Range("D3").Value = "N"
$wb.Worksheets.Item("OpenAccountTest").Range("D4").Value = "Y"
$wb.Worksheets.Item("OpenAccountTest").Range("D5").Value = "Y"

# --- Step 4: cosmetics - column width for the new sheet's first column ---
$wsSuite.Columns.Item(1).ColumnWidth = 21.0221354166667

# --- Step 5: selections (last one selected becomes the active tab) ---
$wb.Worksheets.Item("AddCustomerTest").Range("E1:E5").Select()
$wb.Worksheets.Item("OpenAccountTest").Range("C5").Select()
$wsSuite.Range("B6").Select()
